# Update the approval date on the front matter from "18th February 2020"
# to "28th March 2020", and move the "_GoBack" bookmark (which tracks the
# most recent edit point) from the "Autumn Term 2022" line to sit right
# after the new "March" text, matching where Word leaves it after a live
# edit in that spot.

$d = $word.ActiveDocument

# Locate the date line without disturbing anything else in the document.
$dateRange = $d.Content
$found = $dateRange.Find.Execute("18th February 2020", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $base = $dateRange.Start

    # "18" -> "28"  (keeps the existing run formatting; "th" stays superscript)
    $dayRange = $d.Range($base, $base + 2)
    $dayRange.Text = "28"

    # "February" -> "March" (runs: "28" / "th" (superscript) / " February 2020")
    # "th" occupies 2 chars right after "28", then a space, so "February"
    # starts 5 chars after $base and is 8 characters long.
    $monthRange = $d.Range($base + 5, $base + 13)
    $monthRange.Text = "March"

    # Find exactly where "March" ends now, so the bookmark lands right
    # after it (before " 2020"), regardless of any length changes above.
    $afterMonth = $d.Content
    $afterMonth.Find.Execute("March", $false, $false, $false, $false, `
        $false, $true, 1, $false, "", 0) | Out-Null
    $bookmarkPos = $afterMonth.End

    # The "_GoBack" bookmark currently sits at the end of the "Autumn Term
    # 2022" paragraph (Word's record of the last edit location before this
    # one) - drop it there and re-create it at the new edit location.
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }

    $newBookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
    $d.Bookmarks.Add("_GoBack", $newBookmarkRange)
}
